$d = $word.ActiveDocument

# 1) Update the letter date: "September 19, 2025" -> "September 21, 2025"
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -like "*September 19, 2025*") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# 2) Split the one-line mailing address into two paragraphs. The address text
#    appears twice in the document (once in the letterhead block, once inside
#    the "PROPERTY ADDRESS" table) -- only the non-table occurrence changes.
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -like "*979 Story Road, San Jose CA 95122*" -and $p.Range.Information(12) -eq $false) {
        $p.Range.Text = "979 Story Road"
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.Text = "San Jose, CA 95122"
        break
    }
}

# 3) Remove the now-redundant empty "NoSpacing" paragraph that immediately
#    follows the "Board of Directors" line.
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -like "*Board of Directors*") {
        $next = $p.Next()
        if ($next -ne $null) {
            $next.Range.Delete()
        }
        break
    }
}

$d.Save()
